# Insert a new data row at row 760 ("2026/02/02" / "月" / 19 / 201),
# pushing the existing rows 760:801 down to 761:802.
#
# Sheet layout: A=日付 (date, stored as literal text, e.g. "2026/02/02"),
#               B=曜日 (weekday text), C=時刻 (number), D=ランキング (number).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 760 downward by inserting a blank row at 760.
$ws.Rows.Item(760).Insert()

# The new row's date must stay a literal text string like "2026/02/02"
# (matching every other date cell in column A), not get auto-converted
# to a date serial number by Excel's smart input. Force Text format on
# that one cell before writing to it.
$ws.Cells.Item(760, 1).NumberFormat = "@"
$ws.Cells.Item(760, 1).Value = "2026/02/02"
$ws.Cells.Item(760, 2).Value = "月"
$ws.Cells.Item(760, 3).Value = 19
$ws.Cells.Item(760, 4).Value = 201

# Re-sync the new cell's format with its neighbor above (A759) so the new
# row carries the same (default/general) style as the rest of the sheet,
# rather than leaving behind the one-off "@" text format we applied above.
$ws.Cells.Item(759, 1).Copy()
$ws.Cells.Item(760, 1).PasteSpecial(-4122)
